$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NIK value for row 2 (was "EN-4-070")
$ws.Cells.Replace("EN-4-070", "KJ-6-168")

# Update Nama Karyawan value for row 2 (was "Retno")
$ws.Cells.Replace("Retno", "Solihin")

# Move the active selection/cursor to E9 (was D13)
$ws.Range("E9").Select()
